# Screw-Head_Mount_Precise-Control_Knobs.xlsx -- "working phase" update.
#
# For both data sheets (Table_1, Table_2) the table gains two new rows at
# the very top:
#   Row 1  -> a machine-readable numeric column index (0, 1, 2, ...) using
#             the same bold/bordered/centered style the old header row had.
#   Row 2  -> a new sub-header row that tags a couple of columns
#             ("Head" / "Stud") and is blank everywhere else.
# Everything that used to live in rows 1..N shifts down two rows, and the
# (former) textual header row loses its special header styling once it
# becomes an ordinary data row.

function Update-Table($ws, $lastCol, $rowTags, $dropCols) {
    # Shift the whole table down by two rows.
    $ws.Rows("1:2").Insert()

    # The old header-label row is now row 3 - it's a plain row now, so
    # strip the bold/border/alignment formatting it used to carry.
    $oldHeaderRow = $ws.Range($ws.Cells.Item(3, 1), $ws.Cells.Item(3, $lastCol))
    $oldHeaderRow.ClearFormats()

    # The old header row also carried a couple of trailing metadata columns
    # (e.g. "thread_size" / "material_surface") that don't belong on row 3
    # anymore - blank them out, but keep the (now-empty) cell in place,
    # matching its still-present-but-blank neighbours on the row.
    foreach ($col in $dropCols) {
        $cell = $ws.Cells.Item(3, $col)
        $cell.ClearContents()
        $cell.NumberFormat = "General"
        $cell.Style = "Normal"
    }

    # New row 1: numeric column indices, keeping the header look (bold,
    # thin box border, centered/top aligned) that row 1 used to have.
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item(1, $c).Value = $c - 1
    }
    $newHeaderRow = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $lastCol))
    $newHeaderRow.Font.Bold = $true
    $newHeaderRow.Borders.LineStyle = 1
    $newHeaderRow.HorizontalAlignment = -4108
    $newHeaderRow.VerticalAlignment = -4160

    # New row 2: blank sub-header row, with a couple of tag cells. Touch
    # every cell so the row is fully materialized (matches the rest of the
    # sheet, where every column has an explicit, if empty, cell), then
    # reset back to the default/unstyled look.
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item(2, $c).NumberFormat = "General"
    }
    $newTagRow = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $newTagRow.Style = "Normal"

    foreach ($col in $rowTags.Keys) {
        $ws.Cells.Item(2, $col).Value = $rowTags[$col]
    }
}

$wb = $excel.ActiveWorkbook

# Table_1: A..M (13 cols) -> "Head" in col B(2), "Stud" in col E(5);
# drop the trailing "thread_size"/"material_surface" metadata in L(12)/M(13)
$ws1 = $wb.Worksheets.Item("Table_1")
Update-Table $ws1 13 @{ 2 = "Head"; 5 = "Stud" } @(12, 13)

# Table_2: A..K (11 cols) -> "Head" in col B(2);
# drop the trailing "thread_size"/"material_surface" metadata in J(10)/K(11)
$ws2 = $wb.Worksheets.Item("Table_2")
Update-Table $ws2 11 @{ 2 = "Head" } @(10, 11)
